$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param($ws, [string]$cellRef, $value)
    if ($null -eq $value) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $value
    }
}

$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H40" = 1499.5
    "I40" = 1499.5
    "J40" = 0
    "K40" = 1499.5
    "L40" = 0
    "M40" = -1324.5
    "N40" = $null
    "H43" = 54172744
    "J43" = 9359.200000000001
    "L43" = 9359.200000000001
    "N43" = -9497.200000000001
    "H64" = 7212.0586
    "J64" = 9150.637000000001
    "L64" = 9150.637000000001
    "N64" = -9646.637000000001
    "H67" = 7212.0586
    "J67" = 9150.637000000001
    "L67" = 9150.637000000001
    "N67" = -10866.637
    "H103" = 505.6875
    "J103" = 1200
    "L103" = 3600
    "N103" = -4772
    "H107" = 660.8461
    "I107" = 711.2353000000001
    "J107" = 565.6667
    "K107" = 711.2353000000001
    "L107" = 565.6667
    "M107" = 1208.7647
    "N107" = -4405.6667
    "H112" = 3680
    "J112" = 3838.7097
    "L112" = 11516.1291
    "N112" = -13732.1291
    "H113" = 4286.4287
    "J113" = 4000
    "L113" = 4000
    "N113" = -10508
}
foreach ($key in $updates.Keys) {
    Set-CellValue -ws $ws -cellRef $key -value $updates[$key]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H5" = 432.0435
    "I5" = 495.15384
    "K5" = 495.15384
    "M5" = -383.15384
    "H97" = 11863.333
    "I97" = 11863.333
    "K97" = 11863.333
    "M97" = -11367.333
}
foreach ($key in $updates.Keys) {
    Set-CellValue -ws $ws -cellRef $key -value $updates[$key]
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H4" = 432.0435
    "I4" = 495.15384
    "K4" = 495.15384
    "M4" = -380.15384
    "H35" = 60000
    "J35" = 60000
    "L35" = 60000
    "N35" = -60620
    "H105" = 6228.1665
    "I105" = 7465.0625
    "K105" = 7465.0625
    "M105" = -5718.0625
}
foreach ($key in $updates.Keys) {
    Set-CellValue -ws $ws -cellRef $key -value $updates[$key]
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H16" = 2860.3572
    "I16" = 2120.1
    "K16" = 2120.1
    "M16" = -1833.1
    "H22" = 1666.3334
    "J22" = 999
    "L22" = 999
    "N22" = -1699
    "H62" = 2000
    "I62" = 2000
    "J62" = 2000
    "K62" = 2000
    "L62" = 2000
    "M62" = -1376
    "N62" = -3248
    "H65" = 2000
    "I65" = 2000
    "J65" = 2000
    "K65" = 10000
    "L65" = 10000
    "M65" = -6880
    "N65" = -16240
    "H94" = 1644.3
    "J94" = 1367.875
    "L94" = 1367.875
    "N94" = -2269.875
    "H113" = 2860.3572
    "I113" = 2120.1
    "K113" = 2120.1
    "M113" = 49.90000000000009
    "H122" = 1710.12
    "I122" = 2083.2778
    "K122" = 6249.8334
    "M122" = -3799.8334
}
foreach ($key in $updates.Keys) {
    Set-CellValue -ws $ws -cellRef $key -value $updates[$key]
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H82" = 13477
    "J82" = 20995
    "L82" = 62985
    "N82" = -63797
    "H85" = 13477
    "J85" = 20995
    "L85" = 62985
    "N85" = -65793
    "H107" = 76923220
    "I107" = 145.66667
    "K107" = 437.00001
    "M107" = 1482.99999
    "H109" = 14696.12
    "I109" = 1337.8572
    "K109" = 4013.5716
    "M109" = -2973.5716
    "H132" = 1642.7142
    "I132" = 1166.6666
    "J132" = 1999.75
    "K132" = 10499.9994
    "L132" = 17997.75
    "M132" = -7969.999400000001
    "N132" = -23057.75
    "H137" = 4560.6665
    "I137" = 4505.75
    "K137" = 13517.25
    "M137" = -8417.25
}
foreach ($key in $updates.Keys) {
    Set-CellValue -ws $ws -cellRef $key -value $updates[$key]
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H2" = 134.625
    "I2" = 116
    "K2" = 116
    "M2" = -3
    "H80" = 135863.19
    "I80" = 185416
    "K80" = 185416
    "M80" = -184418
    "H83" = 135863.19
    "I83" = 185416
    "K83" = 927080
    "M83" = -922088
    "H95" = 34749.75
    "J95" = 34749.75
    "L95" = 34749.75
    "N95" = -40241.75
    "H122" = 2407.9167
    "I122" = 2209.6
    "K122" = 6628.799999999999
    "M122" = -4178.799999999999
}
foreach ($key in $updates.Keys) {
    Set-CellValue -ws $ws -cellRef $key -value $updates[$key]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H46" = 9183.75
    "I46" = 12698
    "J46" = 5669.5
    "K46" = 12698
    "L46" = 5669.5
    "M46" = -12510
    "N46" = -6045.5
    "H68" = 2377.4285
    "I68" = 2377.4285
    "K68" = 2377.4285
    "M68" = -1628.4285
    "H71" = 2377.4285
    "I71" = 2377.4285
    "K71" = 11887.1425
    "M71" = -8143.1425
    "H132" = 1925.8
    "I132" = 1925.8
    "J132" = 0
    "K132" = 5777.4
    "L132" = 0
    "M132" = -3247.4
    "N132" = $null
}
foreach ($key in $updates.Keys) {
    Set-CellValue -ws $ws -cellRef $key -value $updates[$key]
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H45" = 29894.5
    "I45" = 0
    "J45" = 29894.5
    "K45" = 0
    "L45" = 29894.5
    "M45" = $null
    "N45" = -30876.5
    "H62" = 4251.5
    "I62" = 4001.5
    "K62" = 4001.5
    "M62" = -3377.5
    "H65" = 4251.5
    "I65" = 4001.5
    "K65" = 20007.5
    "M65" = -16887.5
    "H81" = 64404.938
    "J81" = 250646.25
    "L81" = 501292.5
    "N81" = -503414.5
    "H84" = 64404.938
    "J84" = 250646.25
    "L84" = 2506462.5
    "N84" = -2517070.5
    "H96" = 1400
    "J96" = 0
    "L96" = 0
    "N96" = $null
    "H104" = 36963
    "J104" = 36963
    "L104" = 36963
    "N104" = -43951
    "H107" = 772702.4
    "I107" = 496.6129
    "K107" = 1489.8387
    "M107" = 430.1613
    "H132" = 2755.3057
    "I132" = 1684.5294
    "K132" = 5053.5882
    "M132" = -2523.5882
}
foreach ($key in $updates.Keys) {
    Set-CellValue -ws $ws -cellRef $key -value $updates[$key]
}
